$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Excel_File_Name" values entered for banks whose branch data has now
# --- been scraped, in the chronological order they were actually finished
# --- (by the date embedded in each file name) ---
$ws.Range("D28").Value = "mehreiranian_bank_branches_20241105"
$ws.Range("D15").Value = "toseesaderat_bank_branches_20241113"
$ws.Range("D25").Value = "khavarmianeh_bank_branches_20241115"
$ws.Range("D24").Value = "iranzamin_bank_branches_20241117"
$ws.Range("D21").Value = "sina_bank_branches_20241117"
$ws.Range("D6").Value  = "saman_bank_branches_20241119"
$ws.Range("D17").Value = "keshavarzi_bank_branches_20241120"

# --- Highlight the Saderat Bank row (row 9) in yellow as newly in progress ---
$ws.Range("A9:E9").Interior.Color = 65535

# --- Clear the (no longer needed) "no fill" marker on rows that are now
# --- resolved / no longer need tracking highlight ---
$ws.Range("A16:E16").Interior.Pattern = -4142
$ws.Range("A17:E17").Interior.Pattern = -4142
$ws.Range("A21:E21").Interior.Pattern = -4142
$ws.Range("A23:E23").Interior.Pattern = -4142
$ws.Range("A25:E25").Interior.Pattern = -4142
$ws.Range("A27:E27").Interior.Pattern = -4142
$ws.Range("A28:E28").Interior.Pattern = -4142

# --- Update the view: scroll back to the top and move the active selection ---
$ws.Range("D13").Select()
